$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "63.709.73"
$ws.Range("E2").Value = "  +1.27%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.659.77"
$ws.Range("E3").Value = "  +2.94%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
Set-TextValue $ws.Range("D5") "593.76"
$ws.Range("E5").Value = "  +1.72%  "

# Row 6
Set-TextValue $ws.Range("D6") "147.17"
$ws.Range("E6").Value = "  +0.34%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.593"
$ws.Range("E8").Value = "  -0.69%  "

# Row 9
$ws.Range("E9").Value = "  +0.18%  "

# Row 10
Set-TextValue $ws.Range("D10") "5.68"
$ws.Range("E10").Value = "  +0.15%  "

# Row 11
$ws.Range("E11").Value = "  +0.11%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.356"
$ws.Range("E12").Value = "  +0.84%  "

# Row 13
Set-TextValue $ws.Range("D13") "27.78"
$ws.Range("E13").Value = "  +1.85%  "

# Row 14
Set-TextValue $ws.Range("D14") "3.132.89"
$ws.Range("E14").Value = "  +2.79%  "

# Row 15
Set-TextValue $ws.Range("D15") "63.456.00"
$ws.Range("E15").Value = "  +1.00%  "

# Row 16
$ws.Range("E16").Value = "  +0.58%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.640.02"
$ws.Range("E17").Value = "  +1.91%  "

# Row 18
Set-TextValue $ws.Range("D18") "11.43"
$ws.Range("E18").Value = "  +0.97%  "

# Row 19
Set-TextValue $ws.Range("D19") "343.49"
$ws.Range("E19").Value = "  +0.29%  "

# Row 20
Set-TextValue $ws.Range("D20") "4.38"
$ws.Range("E20").Value = "  -0.08%  "

# Row 21
$ws.Range("E21").Value = "  +1.70%  "

# Row 22
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
Set-TextValue $ws.Range("D23") "68.02"
$ws.Range("E23").Value = "  +1.29%  "

# Row 24
Set-TextValue $ws.Range("D24") "1.70"
$ws.Range("E24").Value = "  +6.92%  "

# Row 25
Set-TextValue $ws.Range("D25") "1.60"
$ws.Range("E25").Value = "  +10.95%  "

# Row 26
Set-TextValue $ws.Range("D26") "562.97"
$ws.Range("E26").Value = "  +21.25%  "

# Row 27
$ws.Range("E27").Value = "  -0.22%  "

# Row 28
Set-TextValue $ws.Range("D28") "8.60"
$ws.Range("E28").Value = "  +3.48%  "

# Row 29
$ws.Range("E29").Value = "  +0.29%  "

# Row 30
$ws.Range("E30").Value = "  +1.46%  "

# Row 31
Set-TextValue $ws.Range("D31") "2.00"
$ws.Range("E31").Value = "  +3.70%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.81"
$ws.Range("E32").Value = "  +12.63%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0₃0822"
$ws.Range("E33").Value = "  +0.24%  "

# Row 34
Set-TextValue $ws.Range("D34") "175.50"
$ws.Range("E34").Value = "  +0.29%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D35") "4.93"
$ws.Range("E35").Value = "  +9.27%  "

# Row 36
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D36") "1.00"
$ws.Range("E36").Value = "  -0.05%  "

# Row 37
$ws.Range("E37").Value = "  +0.30%  "

# Row 38
Set-TextValue $ws.Range("D38") "19.23"
$ws.Range("E38").Value = "  +1.04%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.79"
$ws.Range("E39").Value = "  +5.04%  "

# Row 40
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D40") "0.999"
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D41") "170.26"
$ws.Range("E41").Value = "  +7.37%  "

# Row 42
Set-TextValue $ws.Range("D42") "40.49"
$ws.Range("E42").Value = "  +2.92%  "

# Row 43
Set-TextValue $ws.Range("D43") "3.78"
$ws.Range("E43").Value = "  +0.31%  "

# Row 44
Set-TextValue $ws.Range("D44") "22.03"
$ws.Range("E44").Value = "  +4.12%  "

# Row 45
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D45") "0.0558"
$ws.Range("E45").Value = "  +3.12%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D46") "0.632"
$ws.Range("E46").Value = "  -1.01%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.0242"
$ws.Range("E47").Value = "  +2.33%  "

# Row 48
$ws.Range("E48").Value = "  -0.63%  "

# Row 49
Set-TextValue $ws.Range("D49") "18.85"
$ws.Range("E49").Value = "  +2.13%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.74"
$ws.Range("E50").Value = "  +1.26%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D51") "11.35"
$ws.Range("E51").Value = "  -0.58%  "
